$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits right after
#    "item1" in the table of contents. It gets relocated below, and doing
#    the removal first avoids ever having two bookmarks sharing that name
#    at the same time.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Wrap the "yyyy<CJK>MM<CJK>dd<CJK>" date-placeholder paragraph in
#    parentheses, and re-create the "_GoBack" bookmark right after the new
#    closing parenthesis.
# ---------------------------------------------------------------------------

# Run properties shared by every run in that paragraph.
$rPrHint   = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="微軟正黑體" w:hAnsiTheme="minorHAnsi" w:cs="Arial" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'
$rPrPlain  = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="微軟正黑體" w:hAnsiTheme="minorHAnsi" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

# Locate the whole date paragraph ("yyyy<CJK>MM<CJK>dd<CJK>") by finding
# "yyyy" and taking its enclosing paragraph range (text only, excluding the
# paragraph mark).
$dateFind = $d.Content.Duplicate
$dateFind.Find.Execute("yyyy") | Out-Null
$dateParaRange = $dateFind.Duplicate
$dateParaRange.Expand(4) | Out-Null
$pStart = $dateParaRange.Start
$pEnd = $dateParaRange.End - 1

$target = $d.Range($pStart, $pEnd)

# Rebuild the whole paragraph content in one go: "(" + yyyy + year-char +
# MM + month-char + dd + day-char + ")" + relocated bookmark. Replacing the
# complete paragraph body (rather than a sub-range) keeps run ordering
# correct and avoids splitting the paragraph.
$dateXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
    '<w:r>' + $rPrHint + '<w:t>(</w:t></w:r>' + `
    '<w:r>' + $rPrHint + '<w:t>yyyy</w:t></w:r>' + `
    '<w:r w:rsidR="001E60DC" w:rsidRPr="00E00ED6">' + $rPrPlain + '<w:t>年</w:t></w:r>' + `
    '<w:r>' + $rPrPlain + '<w:t>MM</w:t></w:r>' + `
    '<w:r w:rsidR="001E60DC" w:rsidRPr="00E00ED6">' + $rPrPlain + '<w:t>月</w:t></w:r>' + `
    '<w:r>' + $rPrPlain + '<w:t>dd</w:t></w:r>' + `
    '<w:r w:rsidR="001E60DC" w:rsidRPr="00E00ED6">' + $rPrPlain + '<w:t>日</w:t></w:r>' + `
    '<w:r>' + $rPrHint + '<w:t>)</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($dateXml)
